$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.084.11'
$ws.Range("E2").Value = '  -0.05%  '
$ws.Range("D3").Value = '1.790.07'
$ws.Range("E3").Value = '  -0.25%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '''227.43'
$ws.Range("E5").Value = '  +1.71%  '
$ws.Range("D6").Value = '''0.546'
$ws.Range("E6").Value = '  -1.07%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '''32.24'
$ws.Range("E8").Value = '  -0.57%  '
$ws.Range("E9").Value = '  +3.95%  '
$ws.Range("D10").Value = '''0.0687'
$ws.Range("E10").Value = '  -3.43%  '
$ws.Range("E11").Value = '  +1.13%  '
$ws.Range("E12").Value = '  -0.13%  '
$ws.Range("D13").Value = '''11.47'
$ws.Range("E13").Value = '  +4.26%  '
$ws.Range("D14").Value = '1.789.77'
$ws.Range("E14").Value = '  -0.17%  '
$ws.Range("E15").Value = '  -0.18%  '
$ws.Range("D16").Value = '34.072.66'
$ws.Range("E16").Value = '  -0.12%  '
$ws.Range("E17").Value = '  +0.44%  '
$ws.Range("D18").Value = '''68.14'
$ws.Range("E18").Value = '  +0.11%  '
$ws.Range("D19").Value = '''243.67'
$ws.Range("E19").Value = '  -0.22%  '
$ws.Range("D20").Value = '0.0₃0776'
$ws.Range("E20").Value = '  -1.10%  '
$ws.Range("E21").Value = '  +0.01%  '
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("E23").Value = '  +0.54%  '
$ws.Range("E24").Value = '  -2.81%  '
$ws.Range("D25").Value = '''161.85'
$ws.Range("E26").Value = '  +1.80%  '
$ws.Range("D27").Value = '''16.26'
$ws.Range("E27").Value = '  -0.36%  '
$ws.Range("E28").Value = '  +1.16%  '
$ws.Range("E29").Value = '  +0.12%  '
$ws.Range("E30").Value = '  +2.90%  '
$ws.Range("D31").Value = '''0.0517'
$ws.Range("E31").Value = '  -0.55%  '
$ws.Range("E32").Value = '  -0.22%  '
$ws.Range("E33").Value = '  +4.04%  '
$ws.Range("E34").Value = '  +0.83%  '
$ws.Range("D35").Value = '1.408.40'
$ws.Range("E35").Value = '  +1.57%  '
$ws.Range("D36").Value = '''0.656'
$ws.Range("E36").Value = '  +0.74%  '
$ws.Range("E37").Value = '  -0.41%  '
$ws.Range("E38").Value = '  +2.20%  '
$ws.Range("D39").Value = '''2.35'
$ws.Range("E39").Value = '  +7.97%  '
$ws.Range("D40").Value = '''80.15'
$ws.Range("E40").Value = '  +0.48%  '
$ws.Range("E41").Value = '  +0.64%  '
$ws.Range("D42").Value = '''0.925'
$ws.Range("E42").Value = '  +0.55%  '
$ws.Range("E43").Value = '  -0.86%  '
$ws.Range("D44").Value = '''13.44'
$ws.Range("E44").Value = '  +11.43%  '
$ws.Range("E45").Value = '  +1.41%  '
$ws.Range("D46").Value = '''6.06'
$ws.Range("E46").Value = '  +3.42%  '
$ws.Range("E47").Value = '  +1.77%  '
$ws.Range("E48").Value = '  +2.27%  '
$ws.Range("D49").Value = '''107.67'
$ws.Range("E49").Value = '  +0.17%  '
$ws.Range("D50").Value = '1.949.74'
$ws.Range("E50").Value = '  -0.44%  '
$ws.Range("E51").Value = '  +0.13%  '

Write-Host "Update complete"